# Auto-generated Excel COM-interop script to apply Goblin_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 200000260
$ws.Range("I2").Value = 323.75
$ws.Range("K2").Value = 323.75
$ws.Range("M2").Value = -210.75
# Row 33
$ws.Range("H33").Value = 2352.5
$ws.Range("I33").Value = 378
$ws.Range("J33").Value = 4327
$ws.Range("K33").Value = 378
$ws.Range("L33").Value = 4327
$ws.Range("M33").Value = -149
$ws.Range("N33").Value = -4785
# Row 53
$ws.Range("H53").Value = 269.41177
$ws.Range("I53").Value = 107.833336
$ws.Range("J53").Value = 357.54544
$ws.Range("K53").Value = 107.833336
$ws.Range("L53").Value = 357.54544
$ws.Range("M53").Value = 529.166664
$ws.Range("N53").Value = -1631.54544
# Row 62
$ws.Range("H62").Value = 5863.9375
$ws.Range("I62").Value = 1935.625
$ws.Range("K62").Value = 1935.625
$ws.Range("M62").Value = -1311.625
# Row 65
$ws.Range("H65").Value = 5863.9375
$ws.Range("I65").Value = 1935.625
$ws.Range("K65").Value = 9678.125
$ws.Range("M65").Value = -6558.125
# Row 88
$ws.Range("H88").Value = 4378
$ws.Range("J88").Value = 4478.4
$ws.Range("L88").Value = 4478.4
$ws.Range("N88").Value = -5290.4
# Row 91
$ws.Range("H91").Value = 4378
$ws.Range("J91").Value = 4478.4
$ws.Range("L91").Value = 4478.4
$ws.Range("N91").Value = -7286.4
# Row 116
$ws.Range("H116").Value = 5750.0527
$ws.Range("I116").Value = 4560.4287
$ws.Range("K116").Value = 4560.4287
$ws.Range("M116").Value = -1118.4287
# Row 132
$ws.Range("H132").Value = 3228182.8
$ws.Range("I132").Value = 2461.5925
$ws.Range("K132").Value = 7384.7775
$ws.Range("M132").Value = -4854.7775

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2485.0908
$ws.Range("I45").Value = 1613.8572
$ws.Range("K45").Value = 1613.8572
$ws.Range("M45").Value = -1236.8572
# Row 74
$ws.Range("H74").Value = 2271.3572
$ws.Range("I74").Value = 1875.4762
$ws.Range("J74").Value = 3459
$ws.Range("K74").Value = 1875.4762
$ws.Range("L74").Value = 3459
$ws.Range("M74").Value = -1001.4762
$ws.Range("N74").Value = -5207
# Row 77
$ws.Range("H77").Value = 2271.3572
$ws.Range("I77").Value = 1875.4762
$ws.Range("J77").Value = 3459
$ws.Range("K77").Value = 9377.381000000001
$ws.Range("L77").Value = 17295
$ws.Range("M77").Value = -5009.381000000001
$ws.Range("N77").Value = -26031
# Row 112
$ws.Range("H112").Value = 34597
$ws.Range("J112").Value = 34597
$ws.Range("L112").Value = 34597
$ws.Range("N112").Value = -37551
# Row 119
$ws.Range("H119").Value = 105750
$ws.Range("J119").Value = 105750
$ws.Range("L119").Value = 105750
$ws.Range("N119").Value = -115426
# Row 132
$ws.Range("H132").Value = 2143.6316
$ws.Range("I132").Value = 2242.1875
$ws.Range("J132").Value = 1618
$ws.Range("K132").Value = 6726.5625
$ws.Range("L132").Value = 4854
$ws.Range("M132").Value = -4196.5625
$ws.Range("N132").Value = -9914

$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 15615.444
$ws.Range("J100").Value = 15615.444
$ws.Range("L100").Value = 15615.444
$ws.Range("N100").Value = -17779.444
# Row 107
$ws.Range("H107").Value = 3284.24
$ws.Range("I107").Value = 1545.0714
$ws.Range("J107").Value = 5497.727
$ws.Range("K107").Value = 1545.0714
$ws.Range("L107").Value = 5497.727
$ws.Range("M107").Value = 374.9286
$ws.Range("N107").Value = -9337.726999999999
# Row 134
$ws.Range("H134").Value = 2087137.4
$ws.Range("I134").Value = 4581
$ws.Range("J134").Value = 6252250
$ws.Range("K134").Value = 13743
$ws.Range("L134").Value = 18756750
$ws.Range("M134").Value = -11208
$ws.Range("N134").Value = -18761820

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 9988
$ws.Range("I17").Value = 9988
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 9988
$ws.Range("M17").Value = -9814
$ws.Range("N17").ClearContents()
# Row 25
$ws.Range("H25").Value = 5560
$ws.Range("J25").Value = 3950
$ws.Range("L25").Value = 3950
$ws.Range("N25").Value = -4298
# Row 82
$ws.Range("H82").Value = 96223.5
$ws.Range("J82").Value = 96223.5
$ws.Range("L82").Value = 96223.5
$ws.Range("N82").Value = -96945.5
# Row 85
$ws.Range("H85").Value = 96223.5
$ws.Range("J85").Value = 96223.5
$ws.Range("L85").Value = 96223.5
$ws.Range("N85").Value = -98719.5
# Row 94
$ws.Range("H94").Value = 1474.3572
$ws.Range("J94").Value = 1852.5555
$ws.Range("L94").Value = 1852.5555
$ws.Range("N94").Value = -2754.5555
# Row 96
$ws.Range("H96").Value = 45208
$ws.Range("J96").Value = 45208
$ws.Range("L96").Value = 45208
$ws.Range("N96").Value = -50700
# Row 99
$ws.Range("H99").Value = 3098144.5
$ws.Range("I99").Value = 1803.1
$ws.Range("J99").Value = 7521489.5
$ws.Range("K99").Value = 1803.1
$ws.Range("L99").Value = 7521489.5
$ws.Range("M99").Value = -305.0999999999999
$ws.Range("N99").Value = -7524485.5
# Row 126
$ws.Range("H126").Value = 3098144.5
$ws.Range("I126").Value = 1803.1
$ws.Range("J126").Value = 7521489.5
$ws.Range("K126").Value = 5409.299999999999
$ws.Range("L126").Value = 22564468.5
$ws.Range("M126").Value = -2939.299999999999
$ws.Range("N126").Value = -22569408.5
# Row 132
$ws.Range("H132").Value = 2244.72
$ws.Range("I132").Value = 1760
$ws.Range("J132").Value = 5799.3335
$ws.Range("K132").Value = 5280
$ws.Range("L132").Value = 17398.0005
$ws.Range("M132").Value = -2750
$ws.Range("N132").Value = -22458.0005
# Row 134
$ws.Range("H134").Value = 1973.2593
$ws.Range("I134").Value = 1638.6364
$ws.Range("K134").Value = 4915.9092
$ws.Range("M134").Value = -2380.9092

$ws = $wb.Worksheets.Item("CUL")
# Row 126
$ws.Range("H126").Value = 1964.5
$ws.Range("I126").Value = 1964.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5893.5
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 153849890
$ws.Range("I80").Value = 222224420
$ws.Range("K80").Value = 222224420
$ws.Range("M80").Value = -222223422
# Row 83
$ws.Range("H83").Value = 153849890
$ws.Range("I83").Value = 222224420
$ws.Range("K83").Value = 1111122100
$ws.Range("M83").Value = -1111117108
# Row 102
$ws.Range("H102").Value = 2383
$ws.Range("I102").Value = 1840.6364
$ws.Range("K102").Value = 1840.6364
$ws.Range("M102").Value = -218.6364000000001
# Row 122
$ws.Range("H122").Value = 1562.1765
$ws.Range("I122").Value = 1247.3334
$ws.Range("K122").Value = 3742.0002
$ws.Range("M122").Value = -1292.0002
# Row 126
$ws.Range("H126").Value = 3292.1
$ws.Range("I126").Value = 1836.8334
$ws.Range("J126").Value = 5475
$ws.Range("K126").Value = 5510.5002
$ws.Range("L126").Value = 16425
$ws.Range("M126").Value = -3040.5002
$ws.Range("N126").Value = -21365
# Row 132
$ws.Range("H132").Value = 31252764
$ws.Range("I132").Value = 34484880
$ws.Range("K132").Value = 103454640
$ws.Range("M132").Value = -103452110

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1901.7693
$ws.Range("I22").Value = 382.57144
$ws.Range("J22").Value = 3674.1667
$ws.Range("K22").Value = 382.57144
$ws.Range("L22").Value = 3674.1667
$ws.Range("M22").Value = -87.57144
$ws.Range("N22").Value = -4264.1667
# Row 27
$ws.Range("H27").Value = 1901.7693
$ws.Range("I27").Value = 382.57144
$ws.Range("J27").Value = 3674.1667
$ws.Range("K27").Value = 382.57144
$ws.Range("L27").Value = 3674.1667
$ws.Range("M27").Value = -275.57144
$ws.Range("N27").Value = -3888.1667
# Row 40
$ws.Range("H40").Value = 5721.643
$ws.Range("I40").Value = 5483
$ws.Range("K40").Value = 5483
$ws.Range("M40").Value = -5347
# Row 93
$ws.Range("H93").Value = 5210.913
$ws.Range("I93").Value = 2085.2
$ws.Range("J93").Value = 7615.3076
$ws.Range("K93").Value = 2085.2
$ws.Range("L93").Value = 7615.3076
$ws.Range("M93").Value = -837.1999999999998
$ws.Range("N93").Value = -10111.3076
# Row 132
$ws.Range("H132").Value = 2703.25
$ws.Range("I132").Value = 2000.25
$ws.Range("J132").Value = 3054.75
$ws.Range("K132").Value = 6000.75
$ws.Range("L132").Value = 9164.25
$ws.Range("M132").Value = -3470.75
$ws.Range("N132").Value = -14224.25

$ws = $wb.Worksheets.Item("WVR")
# Row 112
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
# Row 126
$ws.Range("I126").Value = 1835.6428
$ws.Range("K126").Value = 5506.928400000001
$ws.Range("M126").Value = -3036.928400000001
# Row 132
$ws.Range("H132").Value = 37039948
$ws.Range("I132").Value = 66669908
$ws.Range("J132").Value = 2494.5
$ws.Range("K132").Value = 200009724
$ws.Range("L132").Value = 7483.5
$ws.Range("M132").Value = -200007194
$ws.Range("N132").Value = -12543.5
